$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F holds "dSF" values; repull/push updated data (mean calculation)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -11
